# Add the new "Camberwell" exposure-site row. This is inserted above the
# existing "Cape Schank" row (row 4), pushing it and every row below it
# down by one. A handful of "Exposure period" cells in the pushed-down
# rows also get their date normalised from a 4-digit year ("2020") to a
# 2-digit year ("20") to match the rest of the sheet's formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 4, shifting rows 4-23 down to 5-24.
$ws.Rows("4:4").Insert()

# Populate the new row with the Camberwell entry.
$ws.Range("A4").Value = "Camberwell"
$ws.Range("B4").Value = "Fu Lin Asian Grocery Supermarket  1397 Toorak Road, Camberwell VIC 3124"
$ws.Range("C4").Value = "30/12/20 2:30pm-2:45pm"
$ws.Range("D4").Value = "Case shopped"

# Normalise the 4-digit year to a 2-digit year in the exposure-period text
# for the rows that were pushed down (now rows 11-16 and 21-22).
$ws.Range("C11").Value = "29/12/20 11:15am-12:15pm"
$ws.Range("C12").Value = "30/12/20 5:00pm-6:30pm"
$ws.Range("C13").Value = "30/12/20 11:15am-11:20am"
$ws.Range("C14").Value = "30/12/20 6:00pm-6:15pm"
$ws.Range("C15").Value = "29/12/20 07:30am-08:00am"
$ws.Range("C16").Value = "31/12/20 08:00am-08:30am"
$ws.Range("C21").Value = "30/12/20 2:00pm-2:30pm"
$ws.Range("C22").Value = "30/12/20 11:00am-11:30am"
